$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 2 (shifts everything down by 2)
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Rows.Item(2).EntireRow.Insert()

# New row 2: 049/DR 2
$ws.Range("A2").Value = "049/DR 2"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "BG12456"
$ws.Range("D2").Value = "HAYLALA ONE"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "annuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 100000
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = 190000

# New row 3: 094/DR 1
$ws.Range("A3").Value = "094/DR 1"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BG12456"
$ws.Range("D3").Value = "HAYLALA ONE"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 18500

# Update the totals row (now row 8, previously row 6) to account for the
# newly added rows
$ws.Range("I8").Value = 110000
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 110000
$ws.Range("M8").Value = 264500.12
